$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily-push record ("2026/02/27", "金", 8, 201) is inserted as row 884,
# pushing the existing rows 884-925 down to 885-926.
$ws.Rows.Item(884).Insert()

# Column A holds a literal date-formatted text string (not a real Excel date
# serial). Flip the cell to text first so the "2026/02/27" assignment isn't
# auto-parsed into a date value, then restore the default (no explicit)
# style so it matches the plain, unstyled data cells around it.
$ws.Cells.Item(884, 1).NumberFormat = "@"
$ws.Cells.Item(884, 1).Value = "2026/02/27"
$ws.Cells.Item(884, 1).Style = "Normal"

$ws.Cells.Item(884, 2).Value = "金"
$ws.Cells.Item(884, 3).Value = 8
$ws.Cells.Item(884, 4).Value = 201
